# The sheet is a flat weekly price list for "Mango" at Vega Central Mapocho
# de Santiago. Each row is one price record; columns A,B,C,E,F,G,H,I,J,K,Q,T
# carry the same fixed classification values on every row (market, region,
# product taxonomy, unit), while D (date), L (grade), M..P (prices), R
# (origin country) and S (unit price) vary per record.
#
# This commit adds one new weekly record. It must be inserted right above
# the current row 251 (pushing the existing row 251 and everything after it
# down by one row, growing the sheet from 374 to 375 data+header rows), and
# populated with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 251; rows 251..374 shift down to 252..375.
$ws.Rows(251).Insert()

# Row 252 is the row that used to be row 251 - its fixed classification
# columns are identical on every row of this sheet, so reuse them for the
# newly inserted row.
$srcRow = 252
$newRow = 251
$constCols = @(1, 2, 3, 5, 6, 7, 8, 9, 10, 11, 17, 20)
foreach ($col in $constCols) {
    $ws.Cells.Item($newRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value2
}

# New record's own data: date, grade, low/high/average price, origin, unit price.
$ws.Cells.Item($newRow, 4).Value  = 44636      # D - fecha (2022-03-16)
$ws.Cells.Item($newRow, 12).Value = "Primera"  # L - calidad
$ws.Cells.Item($newRow, 13).Value = 450        # M - volumen / cantidad
$ws.Cells.Item($newRow, 14).Value = 6000       # N - precio minimo
$ws.Cells.Item($newRow, 15).Value = 6000       # O - precio maximo
$ws.Cells.Item($newRow, 16).Value = 6000       # P - precio promedio
$ws.Cells.Item($newRow, 18).Value = "Perú"     # R - origen
$ws.Cells.Item($newRow, 19).Value = 1500       # S - precio unitario
